$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 437, pushing the existing row 437 (and every
# row after it) down by one. Excel's default Insert() shift for an
# EntireRow range is xlShiftDown, which is what we want here.
$ws.Rows.Item(437).Insert()

# Populate the freshly inserted row 437 with a new weekly price-report
# observation for "Feria Lagunitas de Puerto Montt" / Zanahoria. Columns
# that are constant across this market/category's rows (A, B, C, E, F, G,
# H, I, K, L, M, N, P, Q, R) are carried over unchanged from the
# observation that used to sit in row 437 (now shifted to row 438); only
# the date (D), volume (J) and origin (O) differ for the new record.
$ws.Cells.Item(437, 1).Value = 4
$ws.Cells.Item(437, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(437, 3).Value = "Los Lagos"
$ws.Cells.Item(437, 4).Value = 44995
$ws.Cells.Item(437, 5).Value = 10
$ws.Cells.Item(437, 6).Value = 100114013
$ws.Cells.Item(437, 7).Value = "Zanahoria"
$ws.Cells.Item(437, 8).Value = "Sin especificar"
$ws.Cells.Item(437, 9).Value = "Primera"
$ws.Cells.Item(437, 10).Value = 850
$ws.Cells.Item(437, 11).Value = 10000
$ws.Cells.Item(437, 12).Value = 10000
$ws.Cells.Item(437, 13).Value = 10000
$ws.Cells.Item(437, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(437, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(437, 16).Value = 500
$ws.Cells.Item(437, 17).Value = 20
$ws.Cells.Item(437, 18).Value = "Hortaliza"
